$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RawData")

$ws.Range("A2").Value = 111111
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "a"
$ws.Range("D2").Value = "a"
$ws.Range("E2").Value = "Y"
$ws.Range("F2").Value = "N"
$ws.Range("G2").Value = "N"
